$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reproduce the authoring order baked into the shared-strings table:
# A1,B1,C1,E1,F1, A2,B2,C2, D1,D2, E2,F2 (birthdate column was bolted on
# after the first few fields, before the 2nd row's nationality/DNI).
$ws.Range("A1").Value = "Javier"
$ws.Range("B1").Value = "Muhlach"
$ws.Range("C1").Value = "javier@email.com"
$ws.Range("E1").Value = "española"
$ws.Range("F1").Value = "71896514P"

$ws.Range("A2").Value = "Pedro"
$ws.Range("B2").Value = "Perez"
$ws.Range("C2").Value = "pedro@email.com"

# Birthdate column - force plain text storage so "07/01/1995" isn't
# auto-converted into an Excel serial date.
$ws.Range("D1:D2").NumberFormat = "@"
$ws.Range("D1").Value = "07/01/1995"
$ws.Range("D2").Value = "30/03/2004"

$ws.Range("E2").Value = "rumana"
$ws.Range("F2").Value = "61478945J"

# Age column - numeric
$ws.Range("G1").Value = 24
$ws.Range("G2").Value = 5

# E-mail column kept as plain text too, then turned into hyperlinks.
$ws.Range("C1:C2").NumberFormat = "@"
$null = $ws.Hyperlinks.Add($ws.Range("C1"), "mailto:javier@email.com")
$null = $ws.Hyperlinks.Add($ws.Range("C2"), "mailto:pedro@email.com")

# Leave the selection where the author's session ended up
$null = $ws.Range("G6").Select()
